$wb = $excel.ActiveWorkbook

# Row 28 (ALC), G context = 27772
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1535.4
$ws.Range("I28").Value = 1561.5555
$ws.Range("J28").Value = 1468.1428
$ws.Range("K28").Value = 1561.5555
$ws.Range("L28").Value = 1468.1428
$ws.Range("M28").Value = -1076.5555
$ws.Range("N28").Value = -2438.1428

# Row 51 (ALC), G context = 5486
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2764.7058
$ws.Range("J51").Value = 3625
$ws.Range("L51").Value = 3625
$ws.Range("N51").Value = -4593

# Row 62 (ALC), G context = 27781
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 40002904
$ws.Range("J62").Value = 3005.5
$ws.Range("L62").Value = 3005.5
$ws.Range("N62").Value = -4253.5

# Row 65 (ALC), G context = 27781
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 40002904
$ws.Range("J65").Value = 3005.5
$ws.Range("L65").Value = 15027.5
$ws.Range("N65").Value = -21267.5

# Row 86 (ALC), G context = 12603
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5535.625
$ws.Range("I86").Value = 4201
$ws.Range("K86").Value = 4201
$ws.Range("M86").Value = -3078

# Row 89 (ALC), G context = 12603
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 5535.625
$ws.Range("I89").Value = 4201
$ws.Range("K89").Value = 21005
$ws.Range("M89").Value = -15389

# Row 98 (ALC), G context = 36237
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1243.6086
$ws.Range("I98").Value = 1176.0476
$ws.Range("K98").Value = 1176.0476
$ws.Range("M98").Value = 321.9523999999999

# Row 100 (ALC), G context = 19906
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1208.2903
$ws.Range("I100").Value = 1086.7084
$ws.Range("J100").Value = 1625.1428
$ws.Range("K100").Value = 1086.7084
$ws.Range("L100").Value = 1625.1428
$ws.Range("M100").Value = -545.7084
$ws.Range("N100").Value = -2707.1428

# Row 112 (ALC), G context = 27960
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 72987.86
$ws.Range("J112").Value = 113204.555
$ws.Range("L112").Value = 339613.665
$ws.Range("N112").Value = -341829.665

# Row 122 (ALC), G context = 36237
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1243.6086
$ws.Range("I122").Value = 1176.0476
$ws.Range("K122").Value = 3528.142800000001
$ws.Range("M122").Value = -1078.142800000001

# Row 132 (ALC), G context = 44049
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1273.6349
$ws.Range("I132").Value = 1029.1321
$ws.Range("K132").Value = 3087.3963
$ws.Range("M132").Value = -557.3963000000003

# Row 141 (ALC), G context = 44161
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2206
$ws.Range("I141").Value = 1964.3334
$ws.Range("K141").Value = 5893.0002
$ws.Range("M141").Value = -713.0002000000004

# Row 32 (ARM), G context = 44147
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4029.1719
$ws.Range("I32").Value = 3125.1052
$ws.Range("K32").Value = 3125.1052
$ws.Range("M32").Value = -2838.1052

# Row 74 (ARM), G context = 44000
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6143.4595
$ws.Range("I74").Value = 1271.2
$ws.Range("J74").Value = 27024.572
$ws.Range("K74").Value = 1271.2
$ws.Range("L74").Value = 27024.572
$ws.Range("M74").Value = -397.2
$ws.Range("N74").Value = -28772.572

# Row 77 (ARM), G context = 44000
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6143.4595
$ws.Range("I77").Value = 1271.2
$ws.Range("J77").Value = 27024.572
$ws.Range("K77").Value = 6356
$ws.Range("L77").Value = 135122.86
$ws.Range("M77").Value = -1988
$ws.Range("N77").Value = -143858.86

# Row 110 (ARM), G context = 27708
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 9111.177
$ws.Range("I110").Value = 10490.833
$ws.Range("J110").Value = 5800
$ws.Range("K110").Value = 10490.833
$ws.Range("L110").Value = 5800
$ws.Range("M110").Value = -8445.833000000001
$ws.Range("N110").Value = -9890

# Row 132 (ARM), G context = 43997
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3494.25
$ws.Range("I132").Value = 2955.625
$ws.Range("K132").Value = 8866.875
$ws.Range("M132").Value = -6336.875

# Row 36 (BSM), G context = 2320
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 13586.875
$ws.Range("I36").Value = 615.8333
$ws.Range("K36").Value = 615.8333
$ws.Range("M36").Value = -81.83330000000001

# Row 107 (BSM), G context = 27706
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1798.9565
$ws.Range("I107").Value = 1420.0714
$ws.Range("J107").Value = 2388.3333
$ws.Range("K107").Value = 1420.0714
$ws.Range("L107").Value = 2388.3333
$ws.Range("M107").Value = 499.9286
$ws.Range("N107").Value = -6228.3333

# Row 134 (BSM), G context = 43998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1619.9348
$ws.Range("I134").Value = 1634.4773
$ws.Range("K134").Value = 4903.4319
$ws.Range("M134").Value = -2368.4319

# Row 3 (CRP), G context = 3763
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 7087
$ws.Range("I3").Value = 750
$ws.Range("J3").Value = 9199.333000000001
$ws.Range("K3").Value = 750
$ws.Range("L3").Value = 9199.333000000001
$ws.Range("M3").Value = -637
$ws.Range("N3").Value = -9425.333000000001

# Row 31 (CRP), G context = 44023
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37566.133
$ws.Range("I31").Value = 60518.65
$ws.Range("K31").Value = 60518.65
$ws.Range("M31").Value = -60223.65

# Row 34 (CRP), G context = 44023
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 37566.133
$ws.Range("I34").Value = 60518.65
$ws.Range("K34").Value = 60518.65
$ws.Range("M34").Value = -60316.65

# Row 39 (CRP), G context = 1915
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 52000
$ws.Range("I39").Value = 4000
$ws.Range("J39").Value = 100000
$ws.Range("K39").Value = 4000
$ws.Range("L39").Value = 100000
$ws.Range("M39").Value = -3609
$ws.Range("N39").Value = -100782

# Row 49 (CRP), G context = 1915
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 52000
$ws.Range("I49").Value = 4000
$ws.Range("J49").Value = 100000
$ws.Range("K49").Value = 4000
$ws.Range("L49").Value = 100000
$ws.Range("M49").Value = -3818
$ws.Range("N49").Value = -100364

# Row 80 (CRP), G context = 12015
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 40125
$ws.Range("J80").Value = 40125
$ws.Range("L80").Value = 40125
$ws.Range("N80").Value = -42371

# Row 83 (CRP), G context = 12015
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 40125
$ws.Range("J83").Value = 40125
$ws.Range("L83").Value = 120375
$ws.Range("N83").Value = -131607

# Row 86 (CRP), G context = 12584
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7898
$ws.Range("J86").Value = 6847
$ws.Range("L86").Value = 6847
$ws.Range("N86").Value = -9093

# Row 89 (CRP), G context = 12584
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 7898
$ws.Range("J89").Value = 6847
$ws.Range("L89").Value = 34235
$ws.Range("N89").Value = -45467

# Row 107 (CRP), G context = 27689
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1028.2106
$ws.Range("I107").Value = 530.625
$ws.Range("K107").Value = 530.625
$ws.Range("M107").Value = 1389.375

# Row 132 (CRP), G context = 44019
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4125.909
$ws.Range("I132").Value = 3975.3333
$ws.Range("J132").Value = 4803.5
$ws.Range("K132").Value = 11925.9999
$ws.Range("L132").Value = 14410.5
$ws.Range("M132").Value = -9395.999899999999
$ws.Range("N132").Value = -19470.5

# Row 134 (CRP), G context = 44020
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 19069.592
$ws.Range("I134").Value = 7251.722
$ws.Range("J134").Value = 72250
$ws.Range("K134").Value = 21755.166
$ws.Range("L134").Value = 216750
$ws.Range("M134").Value = -19220.166
$ws.Range("N134").Value = -221820

# Row 55 (CUL), G context = 4733
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2384.1667
$ws.Range("J55").Value = 3501.25
$ws.Range("L55").Value = 10503.75
$ws.Range("N55").Value = -10857.75

# Row 129 (CUL), G context = 36054
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 561.3333
$ws.Range("I129").Value = 561.3333
$ws.Range("K129").Value = 1683.9999
$ws.Range("M129").Value = 3316.0001

# Row 21 (GSM), G context = 4430
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 2013750

# Row 30 (GSM), G context = 4430
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 2013750

# Row 122 (GSM), G context = 36182
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3060.0908
$ws.Range("I122").Value = 1749
$ws.Range("J122").Value = 3809.2856
$ws.Range("K122").Value = 5247
$ws.Range("L122").Value = 11427.8568
$ws.Range("M122").Value = -2797
$ws.Range("N122").Value = -16327.8568

# Row 23 (LTW), G context = 4097
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 8336567.5
$ws.Range("I23").Value = 3351
$ws.Range("J23").Value = 25003000
$ws.Range("K23").Value = 3351
$ws.Range("L23").Value = 25003000
$ws.Range("M23").Value = -3121
$ws.Range("N23").Value = -25003460

# Row 40 (LTW), G context = 36248
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6172.533
$ws.Range("I40").Value = 4758.5
$ws.Range("J40").Value = 9000.6
$ws.Range("K40").Value = 4758.5
$ws.Range("L40").Value = 9000.6
$ws.Range("M40").Value = -4622.5
$ws.Range("N40").Value = -9272.6

# Row 87 (LTW), G context = 10926
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()

# Row 90 (LTW), G context = 10926
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()

# Row 100 (LTW), G context = 19995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 900
$ws.Range("I100").Value = 900
$ws.Range("K100").Value = 900
$ws.Range("M100").Value = -359

# Row 122 (LTW), G context = 36247
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4520.2256
$ws.Range("I122").Value = 3665.682
$ws.Range("J122").Value = 6609.1113
$ws.Range("K122").Value = 10997.046
$ws.Range("L122").Value = 19827.3339
$ws.Range("M122").Value = -8547.045999999998
$ws.Range("N122").Value = -24727.3339

# Row 136 (LTW), G context = 44060
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4901.2
$ws.Range("I136").Value = 4069.5908
$ws.Range("J136").Value = 7188.125
$ws.Range("K136").Value = 12208.7724
$ws.Range("L136").Value = 21564.375
$ws.Range("M136").Value = -9658.7724
$ws.Range("N136").Value = -26664.375

# Row 132 (WVR), G context = 44029
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2784.75
$ws.Range("I132").Value = 2751.7
$ws.Range("J132").Value = 2950
$ws.Range("K132").Value = 8255.099999999999
$ws.Range("L132").Value = 8850
$ws.Range("M132").Value = -5725.099999999999
$ws.Range("N132").Value = -13910
